# Update res_bus/vm_pu.xlsx values for "case with 380 kV done"
# Slack bus voltage (column B) changed from 1.05 to 1.02 pu, and the
# resulting bus voltage magnitudes (columns C-F, I-N) were recomputed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030657477353299
$ws.Range("D2").Value = 1.039378536622454
$ws.Range("E2").Value = 1.030329742727043
$ws.Range("F2").Value = 1.044746256146077
$ws.Range("I2").Value = 1.025584116979388
$ws.Range("J2").Value = 1.035797510637205
$ws.Range("K2").Value = 1.04216389383017
$ws.Range("L2").Value = 1.033141068222965
$ws.Range("M2").Value = 1.04751643971628
$ws.Range("N2").Value = 1.037268461909945

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032342601505273
$ws.Range("D3").Value = 1.041087138633802
$ws.Range("E3").Value = 1.031776279443928
$ws.Range("F3").Value = 1.046737252912972
$ws.Range("I3").Value = 1.025806844320412
$ws.Range("J3").Value = 1.037120432005471
$ws.Range("K3").Value = 1.043680335249586
$ws.Range("L3").Value = 1.034394223107119
$ws.Range("M3").Value = 1.049315664048447
$ws.Range("N3").Value = 1.038593261978295

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033422887584438
$ws.Range("D4").Value = 1.042179758354416
$ws.Range("E4").Value = 1.032703939663169
$ws.Range("F4").Value = 1.048004210847248
$ws.Range("I4").Value = 1.025940443487851
$ws.Range("J4").Value = 1.03796612483043
$ws.Range("K4").Value = 1.04464819487953
$ws.Range("L4").Value = 1.035196268436914
$ws.Range("M4").Value = 1.050458192424207
$ws.Range("N4").Value = 1.039440155784068

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033874658156791
$ws.Range("D5").Value = 1.042636034448037
$ws.Range("E5").Value = 1.033091960325584
$ws.Range("F5").Value = 1.048531787766076
$ws.Range("I5").Value = 1.025994100800968
$ws.Range("J5").Value = 1.038319214212855
$ws.Range("K5").Value = 1.045051919226941
$ws.Range("L5").Value = 1.035531363304241
$ws.Range("M5").Value = 1.050933374368599
$ws.Range("N5").Value = 1.039793746593915

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033950373796635
$ws.Range("D6").Value = 1.042712466865067
$ws.Range("E6").Value = 1.033156996169936
$ws.Range("F6").Value = 1.048620075645964
$ws.Range("I6").Value = 1.026002963335454
$ws.Range("J6").Value = 1.03837835730214
$ws.Range("K6").Value = 1.045119521893967
$ws.Range("L6").Value = 1.035587505754758
$ws.Range("M6").Value = 1.051012859831812
$ws.Range("N6").Value = 1.039852973673171

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033428933475103
$ws.Range("D7").Value = 1.042185867112179
$ws.Range("E7").Value = 1.032709132106133
$ws.Range("F7").Value = 1.04801128013035
$ws.Range("I7").Value = 1.025941170299361
$ws.Range("J7").Value = 1.037970852372974
$ws.Range("K7").Value = 1.044653601844401
$ws.Range("L7").Value = 1.035200754145649
$ws.Range("M7").Value = 1.050464561940751
$ws.Range("N7").Value = 1.039444890040265

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031229087275287
$ws.Range("D8").Value = 1.039958675739219
$ws.Range("E8").Value = 1.030820353333613
$ws.Range("F8").Value = 1.045423580101641
$ws.Range("I8").Value = 1.025661573899713
$ws.Range("J8").Value = 1.036246755985618
$ws.Range("K8").Value = 1.042679176858655
$ws.Range("L8").Value = 1.03356642385337
$ws.Range("M8").Value = 1.048129021609352
$ws.Range("N8").Value = 1.037718345238282

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027273573341372
$ws.Range("D9").Value = 1.035932854523206
$ws.Range("E9").Value = 1.027426730962282
$ws.Range("F9").Value = 1.040697386390989
$ws.Range("I9").Value = 1.025087769245154
$ws.Range("J9").Value = 1.033128069238236
$ws.Range("K9").Value = 1.039095651353641
$ws.Range("L9").Value = 1.030617531845526
$ws.Range("M9").Value = 1.043844716449941
$ws.Range("N9").Value = 1.034595229598

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024580825822473
$ws.Range("D10").Value = 1.033178006588578
$ws.Range("E10").Value = 1.025118261881021
$ws.Range("F10").Value = 1.03743056685633
$ws.Range("I10").Value = 1.024649835117643
$ws.Range("J10").Value = 1.030992492110232
$ws.Range("K10").Value = 1.036633718205619
$ws.Range("L10").Value = 1.028603195537088
$ws.Range("M10").Value = 1.040871061958452
$ws.Range("N10").Value = 1.03245661970552

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023401010703701
$ws.Range("D11").Value = 1.031967601005288
$ws.Range("E11").Value = 1.024107243372021
$ws.Range("F11").Value = 1.035987474428421
$ws.Range("I11").Value = 1.02444684781992
$ws.Range("J11").Value = 1.030053835968519
$ws.Range("K11").Value = 1.035549717990071
$ws.Range("L11").Value = 1.02771900262538
$ws.Range("M11").Value = 1.039554624418762
$ws.Range("N11").Value = 1.031516630564421

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022960641259292
$ws.Range("D12").Value = 1.031515304717901
$ws.Range("E12").Value = 1.023729942514836
$ws.Range("F12").Value = 1.03544706679673
$ws.Range("I12").Value = 1.024369421693509
$ws.Range("J12").Value = 1.029703034651682
$ws.Range("K12").Value = 1.035144312695262
$ws.Range("L12").Value = 1.02738873179055
$ws.Range("M12").Value = 1.039061223190734
$ws.Range("N12").Value = 1.031165331069477

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023055199454487
$ws.Range("D13").Value = 1.031612446929807
$ws.Range("E13").Value = 1.023810955388025
$ws.Range("F13").Value = 1.035563185739772
$ws.Range("I13").Value = 1.024386121961065
$ws.Range("J13").Value = 1.029778380492864
$ws.Range("K13").Value = 1.035231399371132
$ws.Range("L13").Value = 1.027459660093152
$ws.Range("M13").Value = 1.039167260788248
$ws.Range("N13").Value = 1.031240783910395

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023364653477785
$ws.Range("D14").Value = 1.031930269461505
$ws.Range("E14").Value = 1.024076091812398
$ws.Range("F14").Value = 1.035942894038076
$ws.Range("I14").Value = 1.02444048923928
$ws.Range("J14").Value = 1.030024882622044
$ws.Range("K14").Value = 1.03551626371023
$ws.Range("L14").Value = 1.027691740175656
$ws.Range("M14").Value = 1.039513930359664
$ws.Range("N14").Value = 1.031487636100872

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023555033917353
$ws.Range("D15").Value = 1.032125731000485
$ws.Range("E15").Value = 1.024239216028688
$ws.Range("F15").Value = 1.0361762617989
$ws.Range("I15").Value = 1.02447371739471
$ws.Range("J15").Value = 1.030176475300965
$ws.Range("K15").Value = 1.035691410543183
$ws.Range("L15").Value = 1.027834486922843
$ws.Range("M15").Value = 1.039726936852581
$ws.Range("N15").Value = 1.031639444058785

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024658830117922
$ws.Range("D16").Value = 1.033257962417113
$ws.Range("E16").Value = 1.025185115075969
$ws.Range("F16").Value = 1.037525731196826
$ws.Range("I16").Value = 1.024663023536947
$ws.Range("J16").Value = 1.031054489882213
$ws.Range("K16").Value = 1.036705276133188
$ws.Range("L16").Value = 1.028661620602394
$ws.Range("M16").Value = 1.040957814940944
$ws.Range("N16").Value = 1.032518705521449

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02534746981081
$ws.Range("D17").Value = 1.033963442396335
$ws.Range("E17").Value = 1.025775359391482
$ws.Range("F17").Value = 1.038364514269456
$ws.Range("I17").Value = 1.024778179663956
$ws.Range("J17").Value = 1.031601480210056
$ws.Range("K17").Value = 1.037336395554062
$ws.Range("L17").Value = 1.029177224374898
$ws.Range("M17").Value = 1.041722134867205
$ws.Range("N17").Value = 1.033066472638298

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025747810383866
$ws.Range("D18").Value = 1.034373249156852
$ws.Range("E18").Value = 1.02611853867428
$ws.Range("F18").Value = 1.03885101516579
$ws.Range("I18").Value = 1.02484406081666
$ws.Range("J18").Value = 1.031919189101565
$ws.Range("K18").Value = 1.037702787505319
$ws.Range("L18").Value = 1.029476815264117
$ws.Range("M18").Value = 1.042165174631264
$ws.Range("N18").Value = 1.033384632712874

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025884091976474
$ws.Range("D19").Value = 1.034512698252076
$ws.Range("E19").Value = 1.026235368664053
$ws.Range("F19").Value = 1.039016436141297
$ws.Range("I19").Value = 1.024866306809827
$ws.Range("J19").Value = 1.032027293553053
$ws.Range("K19").Value = 1.037827426140087
$ws.Range("L19").Value = 1.029578773891835
$ws.Range("M19").Value = 1.042315771598271
$ws.Range("N19").Value = 1.033492890685084

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025273723379613
$ws.Range("D20").Value = 1.033887926090444
$ws.Range("E20").Value = 1.025712145887143
$ws.Range("F20").Value = 1.03827480557281
$ws.Range("I20").Value = 1.024765957814302
$ws.Range("J20").Value = 1.031542932407706
$ws.Range("K20").Value = 1.037268861710297
$ws.Range("L20").Value = 1.02912202444195
$ws.Range("M20").Value = 1.041640418230662
$ws.Range("N20").Value = 1.033007841691353

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023273586329216
$ws.Range("D21").Value = 1.031836753638651
$ws.Range("E21").Value = 1.023998064802178
$ws.Range("F21").Value = 1.035831201042654
$ws.Range("I21").Value = 1.024424535575504
$ws.Range("J21").Value = 1.029952353439156
$ws.Range("K21").Value = 1.035432454930857
$ws.Range("L21").Value = 1.027623449590269
$ws.Range("M21").Value = 1.039411967483225
$ws.Range("N21").Value = 1.031415003918226

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022003647650606
$ws.Range("D22").Value = 1.030531458323279
$ws.Range("E22").Value = 1.022910126469834
$ws.Range("F22").Value = 1.03426942652756
$ws.Range("I22").Value = 1.024198127467092
$ws.Range("J22").Value = 1.028939870823251
$ws.Range("K22").Value = 1.034261836426627
$ws.Range("L22").Value = 1.026670554479142
$ws.Range("M22").Value = 1.037985246083536
$ws.Range("N22").Value = 1.030401083460873

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02267805790813
$ws.Range("D23").Value = 1.03122492452058
$ws.Range("E23").Value = 1.02348784834587
$ws.Range("F23").Value = 1.035099790536897
$ws.Range("I23").Value = 1.0243192710006
$ws.Range("J23").Value = 1.029477801423429
$ws.Range("K23").Value = 1.034883940579812
$ws.Range("L23").Value = 1.027176729801385
$ws.Range("M23").Value = 1.038744035504078
$ws.Range("N23").Value = 1.030939777984205

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025307050328423
$ws.Range("D24").Value = 1.033922053870422
$ws.Range("E24").Value = 1.025740712746188
$ws.Range("F24").Value = 1.038315349554819
$ws.Range("I24").Value = 1.024771484319307
$ws.Range("J24").Value = 1.031569391781517
$ws.Range("K24").Value = 1.037299382687761
$ws.Range("L24").Value = 1.029146970473218
$ws.Range("M24").Value = 1.041677351031598
$ws.Range("N24").Value = 1.03303433864051

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028305797554059
$ws.Range("D25").Value = 1.036985911177601
$ws.Range("E25").Value = 1.028312019201325
$ws.Range("F25").Value = 1.041939348513086
$ws.Range("I25").Value = 1.0252458078025
$ws.Range("J25").Value = 1.033944101689407
$ws.Range("K25").Value = 1.040034718137619
$ws.Range("L25").Value = 1.031388268535982
$ws.Range("M25").Value = 1.044972707749164
$ws.Range("N25").Value = 1.035412420908852

